# daily auto push: 2025-10-07 07:26 UTC
# Append one new data row (row 74) to Sheet1, mirroring the existing
# "date / weekday / hour / ranking" row layout used by rows 2-73.
#
# The date column stores plain text like "2025/09/22" (not a real Excel
# date), so we enter it with a leading apostrophe to force text entry
# (exactly like Excel's own quote-prefix mechanism) and then restore the
# cell's style to Normal so no date number-format sticks to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 74

$ws.Cells.Item($row, 1).Value = "'2025/10/07"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = "火"
$ws.Cells.Item($row, 3).Value = 16
$ws.Cells.Item($row, 4).Value = 96
